$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.489.49"
$ws.Range("E2").Value = "  -1.07%  "

$ws.Range("D3").Value = "'1.911.85"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'239.11"
$ws.Range("E5").Value = "  -1.45%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").Value = "'0.4783"
$ws.Range("E7").Value = "  -2.17%  "

$ws.Range("D8").Value = "'0.2844"

$ws.Range("D9").Value = "'0.06704"
$ws.Range("E9").Value = "  -2.95%  "

$ws.Range("D10").Value = "'18.98"
$ws.Range("E10").Value = "  -1.23%  "

$ws.Range("D11").Value = "'102.63"
$ws.Range("E11").Value = "  -2.52%  "

$ws.Range("D12").Value = "'0.07709"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").Value = "'1.914.56"
$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("D14").Value = "'5.208"
$ws.Range("E14").Value = "  -2.78%  "

$ws.Range("D15").Value = "'0.6711"
$ws.Range("E15").Value = "  -3.86%  "

$ws.Range("D16").Value = "'272.29"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").Value = "'30.515.03"
$ws.Range("E17").Value = "  -1.01%  "

$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").Value = "'0.000007465"
$ws.Range("E19").Value = "  -3.11%  "

$ws.Range("D20").Value = "'12.67"
$ws.Range("E20").Value = "  -3.04%  "

$ws.Range("D21").Value = "'5.424"
$ws.Range("E21").Value = "  -1.81%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").Value = "'6.308"
$ws.Range("E23").Value = "  -3.52%  "

$ws.Range("D24").Value = "'9.409"
$ws.Range("E24").Value = "  -3.16%  "

$ws.Range("D25").Value = "'166.77"
$ws.Range("E25").Value = "  +0.22%  "

$ws.Range("D26").Value = "'19.35"
$ws.Range("E26").Value = "  -1.18%  "

$ws.Range("D27").Value = "'2.065"
$ws.Range("E27").Value = "  -4.53%  "

$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("D29").Value = "'0.1005"
$ws.Range("E29").Value = "  -3.00%  "

$ws.Range("D30").Value = "'4.607"
$ws.Range("E30").Value = "  +0.87%  "

$ws.Range("D31").Value = "'1.513"
$ws.Range("E31").Value = "  -2.69%  "

$ws.Range("D32").Value = "'4.242"
$ws.Range("E32").Value = "  -2.65%  "

$ws.Range("D33").Value = "'0.04727"
$ws.Range("E33").Value = "  -2.68%  "

$ws.Range("D34").Value = "'0.7277"
$ws.Range("E34").Value = "  -3.60%  "

$ws.Range("D35").Value = "'1.111"
$ws.Range("E35").Value = "  -3.67%  "

$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("D37").Value = "'0.01922"
$ws.Range("E37").Value = "  -3.63%  "

$ws.Range("D38").Value = "'2.607"
$ws.Range("E38").Value = "  -2.00%  "

$ws.Range("D39").Value = "'6.262"
$ws.Range("E39").Value = "  -3.34%  "

$ws.Range("D40").Value = "'74.76"
$ws.Range("E40").Value = "  -3.16%  "

$ws.Range("D41").Value = "'1.964"
$ws.Range("E41").Value = "  -5.72%  "

$ws.Range("E42").Value = "  -4.99%  "

$ws.Range("D43").Value = "'104.94"
$ws.Range("E43").Value = "  -2.54%  "

$ws.Range("D44").Value = "'0.4264"
$ws.Range("E44").Value = "  -3.14%  "

$ws.Range("D45").Value = "'0.9990"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").Value = "'7.420"
$ws.Range("E46").Value = "  -4.00%  "

$ws.Range("D47").Value = "'0.1201"
$ws.Range("E47").Value = "  -3.48%  "

$ws.Range("D48").Value = "'917.44"
$ws.Range("E48").Value = "  -6.98%  "

$ws.Range("D49").Value = "'34.80"
$ws.Range("E49").Value = "  -3.42%  "

$ws.Range("D50").Value = "'8.797"
$ws.Range("E50").Value = "  -5.13%  "

$ws.Range("D51").Value = "'0.05767"
$ws.Range("E51").Value = "  +0.45%  "
